$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'277.15"
$ws.Range("D3").Value = "'20.96"
$ws.Range("D4").Value = "'6.222"
$ws.Range("D5").Value = "'0.06190"
$ws.Range("D6").Value = "'3.579"
$ws.Range("D7").Value = "'6.565"
$ws.Range("D8").Value = "'1.489"
$ws.Range("D9").Value = "'0.8225"
$ws.Range("D10").Value = "'0.1622"
$ws.Range("D11").Value = "'0.08237"
$ws.Range("D12").Value = "'0.03501"
$ws.Range("D13").Value = "'0.03098"
$ws.Range("D14").Value = "'0.09137"
$ws.Range("D15").Value = "'3.778"
$ws.Range("D16").Value = "'0.001635"
$ws.Range("D17").Value = "'0.04699"
$ws.Range("D18").Value = "'0.006447"
$ws.Range("D19").Value = "'0.006164"
$ws.Range("D20").Value = "'0.001068"
$ws.Range("D22").Value = "'3.806"
$ws.Range("D24").Value = "'0.01383"
$ws.Range("D26").Value = "'0.1224"
$ws.Range("D28").Value = "'0.0002741"
$ws.Range("D40").Value = "'0.04673"
$ws.Range("D41").Value = "'0.007045"
$ws.Range("D42").Value = "'0.1104"
$ws.Range("D44").Value = "'0.01120"
$ws.Range("D45").Value = "'0.00006377"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'0.8464"
$ws.Range("D48").Value = "'0.002353"
$ws.Range("D49").Value = "'0.00001903"
$ws.Range("D50").Value = "'0.01242"
